$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.Value = "'58.183.06"
$c.Style = "Normal"
$ws.Range("E2").Value = '  -3.15%  '
$c = $ws.Range("D3")
$c.Value = "'2.449.24"
$c.Style = "Normal"
$ws.Range("E3").Value = '  -3.79%  '
$ws.Range("E4").Value = '  -0.13%  '
$c = $ws.Range("D5")
$c.Value = "'526.91"
$c.Style = "Normal"
$ws.Range("E5").Value = '  -2.12%  '
$c = $ws.Range("D6")
$c.Value = "'133.34"
$c.Style = "Normal"
$ws.Range("E6").Value = '  -7.70%  '
$ws.Range("E7").Value = '  +0.16%  '
$c = $ws.Range("D8")
$c.Value = "'0.554"
$c.Style = "Normal"
$ws.Range("E8").Value = '  -3.08%  '
$c = $ws.Range("D9")
$c.Value = "'2.452.70"
$c.Style = "Normal"
$ws.Range("E9").Value = '  -4.40%  '
$ws.Range("B10").Value = 'Dogecoin'
$ws.Range("C10").Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$c = $ws.Range("D10")
$c.Value = "'0.0982"
$c.Style = "Normal"
$ws.Range("E10").Value = '  -3.01%  '
$ws.Range("B11").Value = 'TRON'
$ws.Range("C11").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$c = $ws.Range("D11")
$c.Value = "'0.160"
$c.Style = "Normal"
$ws.Range("E11").Value = '  -0.47%  '
$c = $ws.Range("D12")
$c.Value = "'5.29"
$c.Style = "Normal"
$ws.Range("E12").Value = '  -3.22%  '
$c = $ws.Range("D13")
$c.Value = "'0.341"
$c.Style = "Normal"
$ws.Range("E13").Value = '  -5.73%  '
$c = $ws.Range("D14")
$c.Value = "'2.881.76"
$c.Style = "Normal"
$c = $ws.Range("D15")
$c.Value = "'58.075.36"
$c.Style = "Normal"
$ws.Range("E15").Value = '  -3.25%  '
$c = $ws.Range("D16")
$c.Value = "'22.46"
$c.Style = "Normal"
$ws.Range("E16").Value = '  -7.00%  '
$ws.Range("E17").Value = '  -3.99%  '
$c = $ws.Range("D18")
$c.Value = "'2.455.50"
$c.Style = "Normal"
$ws.Range("E18").Value = '  -2.51%  '
$c = $ws.Range("D19")
$c.Value = "'10.66"
$c.Style = "Normal"
$ws.Range("E19").Value = '  -5.35%  '
$ws.Range("B20").Value = 'Polkadot'
$ws.Range("C20").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$c = $ws.Range("D20")
$c.Value = "'4.17"
$c.Style = "Normal"
$ws.Range("E20").Value = '  -4.00%  '
$ws.Range("B21").Value = 'BitcoinCash'
$ws.Range("C21").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$c = $ws.Range("D21")
$c.Value = "'319.27"
$c.Style = "Normal"
$ws.Range("E21").Value = '  -2.37%  '
$ws.Range("E22").Value = '  -0.24%  '
$c = $ws.Range("D23")
$c.Value = "'5.70"
$c.Style = "Normal"
$ws.Range("E23").Value = '  -4.16%  '
$c = $ws.Range("D24")
$c.Value = "'62.36"
$c.Style = "Normal"
$ws.Range("E24").Value = '  -1.25%  '
$c = $ws.Range("D25")
$c.Value = "'0.405"
$c.Style = "Normal"
$ws.Range("E25").Value = '  -7.04%  '
$ws.Range("E26").Value = '  -1.89%  '
$ws.Range("E27").Value = '  -1.31%  '
$c = $ws.Range("D28")
$c.Value = "'7.43"
$c.Style = "Normal"
$ws.Range("E28").Value = '  -7.22%  '
$ws.Range("E29").Value = '  -9.10%  '
$c = $ws.Range("D30")
$c.Value = "'0.0₃0746"
$c.Style = "Normal"
$ws.Range("E30").Value = '  -6.29%  '
$ws.Range("E31").Value = '  -3.96%  '
$c = $ws.Range("D32")
$c.Value = "'162.90"
$c.Style = "Normal"
$ws.Range("E32").Value = '  -1.35%  '
$ws.Range("E33").Value = '  +0.06%  '
$c = $ws.Range("D34")
$c.Value = "'1.07"
$c.Style = "Normal"
$ws.Range("E34").Value = '  -9.60%  '
$ws.Range("E35").Value = '  -8.68%  '
$c = $ws.Range("D36")
$c.Value = "'18.12"
$c.Style = "Normal"
$ws.Range("E36").Value = '  -3.31%  '
$c = $ws.Range("D37")
$c.Value = "'4.00"
$c.Style = "Normal"
$ws.Range("E37").Value = '  -9.64%  '
$c = $ws.Range("D38")
$c.Value = "'1.53"
$c.Style = "Normal"
$ws.Range("E38").Value = '  -6.23%  '
$c = $ws.Range("D39")
$c.Value = "'36.30"
$c.Style = "Normal"
$ws.Range("E39").Value = '  -1.90%  '
$c = $ws.Range("D40")
$c.Value = "'3.51"
$c.Style = "Normal"
$ws.Range("E40").Value = '  -5.77%  '
$c = $ws.Range("D41")
$c.Value = "'0.787"
$c.Style = "Normal"
$ws.Range("E41").Value = '  -6.13%  '
$c = $ws.Range("D42")
$c.Value = "'0.997"
$c.Style = "Normal"
$ws.Range("E42").Value = '  +0.14%  '
$c = $ws.Range("D43")
$c.Value = "'273.56"
$c.Style = "Normal"
$ws.Range("E43").Value = '  -9.16%  '
$c = $ws.Range("D44")
$c.Value = "'5.02"
$c.Style = "Normal"
$ws.Range("E44").Value = '  -9.97%  '
$c = $ws.Range("D45")
$c.Value = "'10.82"
$c.Style = "Normal"
$ws.Range("E45").Value = '  -0.26%  '
$ws.Range("E46").Value = '  -4.15%  '
$c = $ws.Range("D47")
$c.Value = "'0.0919"
$c.Style = "Normal"
$ws.Range("E47").Value = '  -2.05%  '
$c = $ws.Range("D48")
$c.Value = "'119.56"
$c.Style = "Normal"
$ws.Range("E48").Value = '  -5.92%  '
$c = $ws.Range("D49")
$c.Value = "'0.0503"
$c.Style = "Normal"
$ws.Range("E49").Value = '  -3.29%  '
$c = $ws.Range("D50")
$c.Value = "'0.0216"
$c.Style = "Normal"
$ws.Range("E50").Value = '  -5.76%  '
$c = $ws.Range("D51")
$c.Value = "'16.90"
$c.Style = "Normal"
$ws.Range("E51").Value = '  -7.44%  '
